$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows: OM numbers + Status ---
# Row 2
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 1).Value = 5198569
$ws.Cells.Item(2, 2).Value = "Encerrado!"

# Row 3
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 1).Value = 685601364540
$ws.Cells.Item(3, 2).Value = "Encerrado!"

# Row 4
$ws.Cells.Item(4, 1).Value = 685601364299
$ws.Cells.Item(4, 2).Value = "Encerrado!"

# Row 5
$ws.Cells.Item(5, 1).Value = 685601378525
$ws.Cells.Item(5, 2).Value = "Encerrado!"

# Row 6
$ws.Cells.Item(6, 1).Value = 685601381337
$ws.Cells.Item(6, 2).Value = "Ordem pendente!"

# Row 7
$ws.Cells.Item(7, 1).Value = 685601381334
$ws.Cells.Item(7, 2).Value = "Ordem pendente!"

# Row 8
$ws.Cells.Item(8, 1).Value = 685601387728
$ws.Cells.Item(8, 2).Value = "Encerrado!"

# --- Column widths (best effort match of final look) ---
$ws.Columns.Item(1).ColumnWidth = 11.1666666666667

# --- Selection matches the saved view state ---
$ws.Range("F12").Select()
